$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels: "L-V" (Lunes-Viernes) entries become "L-J" (Lunes-Jueves)
$ws.Range("J1").Value = "HORA DE SALIDA L-J"
$ws.Range("I1").Value = "HORA DE ENTRADA L-J"

# Update the active selection to I2
$ws.Range("I2").Select()
